# Ran code for averaged intensities on spiral schemes
# Append three new rows (17-19) to the "UniformF" sheet, mirroring the
# shape/formatting of the existing data rows: an index in column A, a
# scheme label in column B, and the averaged-intensity value (1) across
# columns C..M.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ Row = 17; Index = 15; Label = "HexGrid-90degTilt5degRes" },
    @{ Row = 18; Index = 16; Label = "HexGrid-90degTilt22p5degRes" },
    @{ Row = 19; Index = 17; Label = "HexGrid-60degTilt5degRes" }
)

foreach ($item in $newRows) {
    $r = $item.Row

    $ws.Cells.Item($r, 1).Value = $item.Index
    $ws.Cells.Item($r, 2).Value = $item.Label

    for ($c = 3; $c -le 13; $c++) {
        $ws.Cells.Item($r, $c).Value = 1
    }
}

# Copy the formatting of the last pre-existing row-index cell (A16) down
# onto the newly added rows' index cells so the bold/centered/bordered
# look used for column A is preserved (columns B..M keep the default
# style, matching every other data row).
$ws.Range("A16").Copy() | Out-Null
$ws.Range("A17:A19").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
